# Insert two new blank "Title and Content" slides at the very beginning of
# the deck (new slides 1 and 2); the four existing slides are pushed down
# to positions 3-6, keeping their content untouched.

$p = $ppt.ActivePresentation

# Reuse the same slide layout ("Title and Content") the existing slides use.
$layout = $p.Slides.Item(1).CustomLayout

$p.Slides.Add(1, 2) | Out-Null
$p.Slides.Add(2, 2) | Out-Null
